$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet currently has data rows 2-174 (row 1 = header). A new match
# (id 6818361, played earlier than all the other matches in this batch, per
# the Date column) needs to be inserted in chronological order at row 169;
# the six rows that used to occupy 169-174 shift down to 170-175.
#
# Row 175 does not exist yet, so first clone the two styled cells (A = bold
# bordered sequence-number style, E = date-formatted style) from row 174 down
# onto row 175 before writing any values, so the new row matches the sheet's
# existing per-column formatting exactly.
$ws.Range("A174").Copy() | Out-Null
$ws.Range("A175").PasteSpecial(-4122) | Out-Null
$ws.Range("E174").Copy() | Out-Null
$ws.Range("E175").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 169
$ws.Range("A169").Value = 167
$ws.Range("B169").Value = 6818361
$ws.Range("C169").Value = "Hungary NB I"
$ws.Range("D169").Value = "Hungary NB I"
$ws.Range("E169").Value = 45396.57291666666
$ws.Range("F169").Value = "Mezokovesd Zsory"
$ws.Range("G169").Value = "Puskas Academy"
$ws.Range("H169").Value = 0
$ws.Range("I169").Value = 4
$ws.Range("J169").Value = "A"
$ws.Range("K169").Value = 5.5
$ws.Range("L169").Value = 3.5
$ws.Range("M169").Value = 1.666
$ws.Range("N169").Value = 6.5
$ws.Range("O169").Value = 3.75
$ws.Range("P169").Value = 1.55
$ws.Range("Q169").Value = 1
$ws.Range("R169").Value = 1.9
$ws.Range("S169").Value = 1.95
$ws.Range("T169").Value = 2.5
$ws.Range("U169").Value = 1.85
$ws.Range("V169").Value = 2
$ws.Range("W169").Value = -1
$ws.Range("X169").Value = -1
$ws.Range("Y169").Value = 0.55
$ws.Range("Z169").Value = -1
$ws.Range("AA169").Value = 0.95
$ws.Range("AB169").Value = 0.8500000000000001
$ws.Range("AC169").Value = -1

# Row 170
$ws.Range("A170").Value = 168
$ws.Range("B170").Value = 6818365
$ws.Range("C170").Value = "Hungary NB I"
$ws.Range("D170").Value = "Hungary NB I"
$ws.Range("E170").Value = 45402.39583333334
$ws.Range("F170").Value = "Puskas Academy"
$ws.Range("G170").Value = "Paksi"
$ws.Range("K170").Value = 2
$ws.Range("L170").Value = 3.4
$ws.Range("M170").Value = 3.3
$ws.Range("N170").Value = 1.909
$ws.Range("O170").Value = 3.5
$ws.Range("P170").Value = 3.6
$ws.Range("Q170").Value = -0.5
$ws.Range("R170").Value = 1.925
$ws.Range("S170").Value = 1.925
$ws.Range("T170").Value = 2.75
$ws.Range("U170").Value = 2
$ws.Range("V170").Value = 1.85
$ws.Range("W170").Value = 0
$ws.Range("X170").Value = 0
$ws.Range("Y170").Value = 0
$ws.Range("Z170").Value = 0
$ws.Range("AA170").Value = 0

# Row 171
$ws.Range("A171").Value = 169
$ws.Range("B171").Value = 6818362
$ws.Range("C171").Value = "Hungary NB I"
$ws.Range("D171").Value = "Hungary NB I"
$ws.Range("E171").Value = 45402.5
$ws.Range("F171").Value = "Ferencvarosi TC"
$ws.Range("G171").Value = "Kisvarda FC"
$ws.Range("K171").Value = 1.25
$ws.Range("L171").Value = 5.5
$ws.Range("M171").Value = 9
$ws.Range("N171").Value = 1.181
$ws.Range("O171").Value = 5.75
$ws.Range("P171").Value = 13
$ws.Range("Q171").Value = -2
$ws.Range("R171").Value = 2.05
$ws.Range("S171").Value = 1.8
$ws.Range("T171").Value = 3
$ws.Range("U171").Value = 1.875
$ws.Range("V171").Value = 1.975
$ws.Range("W171").Value = 0
$ws.Range("X171").Value = 0
$ws.Range("Y171").Value = 0
$ws.Range("Z171").Value = 0
$ws.Range("AA171").Value = 0

# Row 172
$ws.Range("A172").Value = 170
$ws.Range("B172").Value = 6818364
$ws.Range("C172").Value = "Hungary NB I"
$ws.Range("D172").Value = "Hungary NB I"
$ws.Range("E172").Value = 45402.60416666666
$ws.Range("F172").Value = "Diosgyori VTK"
$ws.Range("G172").Value = "Debreceni VSC"
$ws.Range("K172").Value = 2.5
$ws.Range("L172").Value = 3.2
$ws.Range("M172").Value = 2.625
$ws.Range("N172").Value = 2.625
$ws.Range("O172").Value = 3.2
$ws.Range("P172").Value = 2.5
$ws.Range("Q172").Value = 0
$ws.Range("R172").Value = 2.025
$ws.Range("S172").Value = 1.825
$ws.Range("T172").Value = 2.75
$ws.Range("U172").Value = 2
$ws.Range("V172").Value = 1.85
$ws.Range("W172").Value = 0
$ws.Range("X172").Value = 0
$ws.Range("Y172").Value = 0
$ws.Range("Z172").Value = 0
$ws.Range("AA172").Value = 0

# Row 173
$ws.Range("A173").Value = 171
$ws.Range("B173").Value = 6818367
$ws.Range("C173").Value = "Hungary NB I"
$ws.Range("D173").Value = "Hungary NB I"
$ws.Range("E173").Value = 45403.37847222222
$ws.Range("F173").Value = "Kecskemeti TE"
$ws.Range("G173").Value = "Mezokovesd Zsory"
$ws.Range("K173").Value = 1.727
$ws.Range("L173").Value = 3.5
$ws.Range("M173").Value = 4.333
$ws.Range("N173").Value = 1.65
$ws.Range("O173").Value = 3.4
$ws.Range("P173").Value = 5.25
$ws.Range("Q173").Value = -0.75
$ws.Range("R173").Value = 1.9
$ws.Range("S173").Value = 1.95
$ws.Range("T173").Value = 2.25
$ws.Range("U173").Value = 1.975
$ws.Range("V173").Value = 1.875
$ws.Range("W173").Value = 0
$ws.Range("X173").Value = 0
$ws.Range("Y173").Value = 0
$ws.Range("Z173").Value = 0
$ws.Range("AA173").Value = 0

# Row 174
$ws.Range("A174").Value = 172
$ws.Range("B174").Value = 6818366
$ws.Range("C174").Value = "Hungary NB I"
$ws.Range("D174").Value = "Hungary NB I"
$ws.Range("E174").Value = 45403.45833333334
$ws.Range("F174").Value = "MOL Fehervar FC"
$ws.Range("G174").Value = "Zalaegerszegi TE"
$ws.Range("K174").Value = 1.909
$ws.Range("L174").Value = 3.4
$ws.Range("M174").Value = 3.6
$ws.Range("N174").Value = 1.909
$ws.Range("O174").Value = 3.4
$ws.Range("P174").Value = 3.6
$ws.Range("Q174").Value = -0.5
$ws.Range("R174").Value = 1.975
$ws.Range("S174").Value = 1.875
$ws.Range("T174").Value = 2.5
$ws.Range("U174").Value = 1.825
$ws.Range("V174").Value = 2.025
$ws.Range("W174").Value = 0
$ws.Range("X174").Value = 0
$ws.Range("Y174").Value = 0
$ws.Range("Z174").Value = 0
$ws.Range("AA174").Value = 0

# Row 175
$ws.Range("A175").Value = 173
$ws.Range("B175").Value = 6818363
$ws.Range("C175").Value = "Hungary NB I"
$ws.Range("D175").Value = "Hungary NB I"
$ws.Range("E175").Value = 45403.5625
$ws.Range("F175").Value = "Ujpest"
$ws.Range("G175").Value = "MTK Budapest"
$ws.Range("K175").Value = 2
$ws.Range("L175").Value = 3.4
$ws.Range("M175").Value = 3.3
$ws.Range("N175").Value = 2.1
$ws.Range("O175").Value = 3.4
$ws.Range("P175").Value = 3.1
$ws.Range("Q175").Value = -0.25
$ws.Range("R175").Value = 1.9
$ws.Range("S175").Value = 1.95
$ws.Range("T175").Value = 2.75
$ws.Range("U175").Value = 1.825
$ws.Range("V175").Value = 2.025
$ws.Range("W175").Value = 0
$ws.Range("X175").Value = 0
$ws.Range("Y175").Value = 0
$ws.Range("Z175").Value = 0
$ws.Range("AA175").Value = 0
